$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for the rows that were repulled / recalculated.
$ws.Range("F11").Value = -5
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = -1
$ws.Range("F21").Value = -2
$ws.Range("F29").Value = 1
$ws.Range("F32").Value = -3
$ws.Range("F42").Value = -6
$ws.Range("F43").Value = 2
$ws.Range("F44").Value = 3
$ws.Range("F47").Value = 2
$ws.Range("F48").Value = -7
$ws.Range("F49").Value = 2
$ws.Range("F52").Value = 2
$ws.Range("F53").Value = -5
$ws.Range("F55").Value = 3
$ws.Range("F56").Value = -1
